# Auto-generated PowerShell COM-interop script
# Applies the "run 162" optimisation_result.xlsx update to before.xlsx

$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# ---- Schedule sheet: update rows 2-4 ----
$wsSchedule.Range("A2").Value = 46060.0625
$wsSchedule.Range("B2").Value = 46060.22916666666
$wsSchedule.Range("C2").Value = 4
$wsSchedule.Range("D2").Value = 15.12
$wsSchedule.Range("E2").Value = 775.5372592500001
$wsSchedule.Range("F2").Value = 51.29214677579366
$wsSchedule.Range("A3").Value = 46060.29166666666
$wsSchedule.Range("C3").Value = 9.5
$wsSchedule.Range("D3").Value = 35.91
$wsSchedule.Range("E3").Value = 646.735635
$wsSchedule.Range("F3").Value = 18.00990350877193
$wsSchedule.Range("A4").Value = 46060.89583333334
$wsSchedule.Range("B4").Value = 46061.0625
$wsSchedule.Range("C4").Value = 4
$wsSchedule.Range("D4").Value = 15.12
$wsSchedule.Range("E4").Value = 687.70106925
$wsSchedule.Range("F4").Value = 45.48287495039683

# ---- Schedule sheet: add new row 5 ----
$wsSchedule.Range("A5").Value = 46061.16666666666
$wsSchedule.Range("B5").Value = 46061.60416666666
$wsSchedule.Range("C5").Value = 10.5
$wsSchedule.Range("D5").Value = 39.69
$wsSchedule.Range("E5").Value = 1220.94551475
$wsSchedule.Range("F5").Value = 30.762043707483
$wsSchedule.Range("A5:B5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---- Detailed sheet: cell updates ----
$wsDetailed.Range("E2").Value = "OFF"
$wsDetailed.Range("E3").Value = "OFF"
$wsDetailed.Range("E4").Value = "OFF"
$wsDetailed.Range("E13").Value = "OFF"
$wsDetailed.Range("E15").Value = "OFF"
$wsDetailed.Range("B37").Value = 115.05557
$wsDetailed.Range("B38").Value = 163.87544
$wsDetailed.Range("B39").Value = 151.54563
$wsDetailed.Range("C39").Value = "historical"
$wsDetailed.Range("B40").Value = 150.68478
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("B41").Value = 159.07424
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 189.86016
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 240.89
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 121.06852
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("B45").Value = 108.89
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("E45").Value = "ON"
$wsDetailed.Range("B46").Value = 89.27254000000001
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("E46").Value = "ON"
$wsDetailed.Range("B47").Value = 69.30265
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("E47").Value = "ON"
$wsDetailed.Range("B48").Value = 79.95028000000001
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("E48").Value = "ON"
$wsDetailed.Range("B49").Value = 82.06932
$wsDetailed.Range("E49").Value = "ON"
$wsDetailed.Range("B50").Value = 84.79000000000001
$wsDetailed.Range("E50").Value = "ON"
$wsDetailed.Range("B51").Value = 85.26964
$wsDetailed.Range("E51").Value = "ON"
$wsDetailed.Range("B52").Value = 105.79
$wsDetailed.Range("E52").Value = "ON"
$wsDetailed.Range("B53").Value = 84.79000000000001
$wsDetailed.Range("B54").Value = 81.37902
$wsDetailed.Range("B55").Value = 78.00005
$wsDetailed.Range("B56").Value = 78.00005
$wsDetailed.Range("B57").Value = 73.20007
$wsDetailed.Range("B58").Value = 69.70384
$wsDetailed.Range("E58").Value = "ON"
$wsDetailed.Range("B59").Value = 69.38724999999999
$wsDetailed.Range("E59").Value = "ON"
$wsDetailed.Range("B60").Value = 66.08642
$wsDetailed.Range("E60").Value = "ON"
$wsDetailed.Range("B61").Value = 73.20010000000001
$wsDetailed.Range("E61").Value = "ON"
$wsDetailed.Range("B62").Value = 69.46913000000001
$wsDetailed.Range("B63").Value = 69.11879999999999
$wsDetailed.Range("B64").Value = 62.28628
$wsDetailed.Range("B65").Value = 57.06007
$wsDetailed.Range("B66").Value = 56.97996
$wsDetailed.Range("B67").Value = 36.07
$wsDetailed.Range("B68").Value = 45.50642
$wsDetailed.Range("B69").Value = 56.98
$wsDetailed.Range("B70").Value = 57.08
$wsDetailed.Range("B71").Value = 56.98
$wsDetailed.Range("B72").Value = 56.98
$wsDetailed.Range("B73").Value = 56.98
$wsDetailed.Range("B74").Value = 57.08
$wsDetailed.Range("B75").Value = 57.08
$wsDetailed.Range("B76").Value = 61.19489
$wsDetailed.Range("B77").Value = 57.08
$wsDetailed.Range("B78").Value = 59.94865
$wsDetailed.Range("B79").Value = 73.20010000000001
$wsDetailed.Range("E79").Value = "OFF"
$wsDetailed.Range("B80").Value = 108.01
$wsDetailed.Range("E80").Value = "OFF"
$wsDetailed.Range("B81").Value = 105
$wsDetailed.Range("E81").Value = "OFF"
$wsDetailed.Range("B82").Value = 74.36649
$wsDetailed.Range("E82").Value = "OFF"
$wsDetailed.Range("B83").Value = 65.35863999999999
$wsDetailed.Range("E83").Value = "OFF"
$wsDetailed.Range("B84").Value = 73.20010000000001
$wsDetailed.Range("E84").Value = "OFF"
$wsDetailed.Range("B85").Value = 103.38668
$wsDetailed.Range("E85").Value = "OFF"
$wsDetailed.Range("B87").Value = 108.01
$wsDetailed.Range("B88").Value = 119.11621
$wsDetailed.Range("B89").Value = 107.76225
$wsDetailed.Range("B90").Value = 111.41472
$wsDetailed.Range("B91").Value = 108.89
$wsDetailed.Range("B92").Value = 108.01
$wsDetailed.Range("B93").Value = 105.79
$wsDetailed.Range("B94").Value = 84.79000000000001
$wsDetailed.Range("B95").Value = 80.14892
$wsDetailed.Range("B96").Value = 82.68552
$wsDetailed.Range("B97").Value = 84.79000000000001
